# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range('D2').Value = '37.052.45'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '2.057.13'
$ws.Range('E3').Value = '  -2.42%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '248.82'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.661'
$ws.Range('D6').NumberFormat = "General"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '55.42'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +15.85%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '61.24'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +3.12%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.379'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.49%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0791'
$ws.Range('D11').NumberFormat = "General"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.96%  '
$ws.Range('E12').Value = '  +5.88%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.08'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +4.51%  '
$ws.Range('D14').Value = '2.358.93'
$ws.Range('E14').Value = '  -2.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.814'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.23'
$ws.Range('D16').NumberFormat = "General"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.64%  '
$ws.Range('D17').Value = '2.058.92'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '37.000.57'
$ws.Range('E18').Value = '  +1.28%  '
$val_D19 = "{0}.0{1}0947" -f '0', $sub3
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = $val_D19
$ws.Range('D19').NumberFormat = "General"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +13.85%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '72.42'
$ws.Range('D20').NumberFormat = "General"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.22'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +6.99%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.33'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.85'
$ws.Range('D23').NumberFormat = "General"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.41%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.42'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '170.32'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.04'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.06'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -6.55%  '
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.55'
$ws.Range('D31').NumberFormat = "General"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0622'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.04'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +10.08%  '
$ws.Range('E34').Value = '  +7.28%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.28'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.78%  '
$ws.Range('E37').Value = '  -6.00%  '
$ws.Range('E38').Value = '  -5.50%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.34'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.105'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +24.62%  '
$ws.Range('E41').Value = '  +11.97%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0223'
$ws.Range('D42').NumberFormat = "General"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '96.29'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.87%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.78'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.15%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.16'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +51.96%  '
$ws.Range('E47').Value = '  +7.81%  '
$ws.Range('D48').Value = '1.296.11'
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('E50').Value = '  -54.59%  '
$ws.Range('E51').Value = '  -2.47%  '
